$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (Dec-2018 and Sep-2018 quarters),
# shifting the existing quarterly columns D:K to F:M.
$ws.Columns("D:E").Insert(-4161)

# The newly inserted columns don't inherit number formatting from the
# insert itself, so copy formats from the (now shifted) F:G columns,
# which still carry the original D:E formatting (date row style / numeric style).
# Use a bounded range (not whole columns) so we don't materialize the
# entire 1,048,576-row column into the sheet.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("F7").Value2 = 43281
$ws.Range("G7").Value2 = 43190
$ws.Range("H7").Value2 = 43100
$ws.Range("I7").Value2 = 43008
$ws.Range("J7").Value2 = 42916
$ws.Range("K7").Value2 = 42825
$ws.Range("L7").Value2 = 42735
$ws.Range("M7").Value2 = 42643
$ws.Range("D8").Value2 = 991300
$ws.Range("E8").Value2 = 1296000
$ws.Range("F8").Value2 = 1332700
$ws.Range("G8").Value2 = 926600
$ws.Range("H8").Value2 = 964300
$ws.Range("I8").Value2 = 1229600
$ws.Range("J8").Value2 = 1275900
$ws.Range("K8").Value2 = 872100
$ws.Range("L8").Value2 = 913600
$ws.Range("M8").Value2 = 1241200
$ws.Range("D9").Value2 = 741700
$ws.Range("E9").Value2 = 977000
$ws.Range("F9").Value2 = 1012000
$ws.Range("G9").Value2 = 695700
$ws.Range("H9").Value2 = 723400
$ws.Range("I9").Value2 = 933700
$ws.Range("J9").Value2 = 965600
$ws.Range("K9").Value2 = 653500
$ws.Range("L9").Value2 = 685500
$ws.Range("M9").Value2 = 939000
$ws.Range("D10").Value2 = 249600
$ws.Range("E10").Value2 = 319000
$ws.Range("F10").Value2 = 320700
$ws.Range("G10").Value2 = 230900
$ws.Range("H10").Value2 = 240900
$ws.Range("I10").Value2 = 295900
$ws.Range("J10").Value2 = 310300
$ws.Range("K10").Value2 = 218600
$ws.Range("L10").Value2 = 228100
$ws.Range("M10").Value2 = 302200
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("F12").Value2 = "NA"
$ws.Range("G12").Value2 = "NA"
$ws.Range("H12").Value2 = "NA"
$ws.Range("I12").Value2 = "NA"
$ws.Range("J12").Value2 = "NA"
$ws.Range("K12").Value2 = "NA"
$ws.Range("L12").Value2 = "NA"
$ws.Range("M12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("F13").Value2 = 0
$ws.Range("G13").Value2 = 0
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("M13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("F14").Value2 = 0
$ws.Range("G14").Value2 = 0
$ws.Range("H14").Value2 = 0
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 0
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 0
$ws.Range("M14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("F15").Value2 = 0
$ws.Range("G15").Value2 = 0
$ws.Range("H15").Value2 = 0
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 0
$ws.Range("K15").Value2 = 0
$ws.Range("L15").Value2 = 0
$ws.Range("M15").Value2 = 0
$ws.Range("D17").Value2 = 932800
$ws.Range("E17").Value2 = 1173700
$ws.Range("F17").Value2 = 1195400
$ws.Range("G17").Value2 = 872600
$ws.Range("H17").Value2 = 903000
$ws.Range("I17").Value2 = 1115100
$ws.Range("J17").Value2 = 1146600
$ws.Range("K17").Value2 = 823400
$ws.Range("L17").Value2 = 855500
$ws.Range("M17").Value2 = 1121900
$ws.Range("D18").Value2 = 58500
$ws.Range("E18").Value2 = 122300
$ws.Range("F18").Value2 = 137300
$ws.Range("G18").Value2 = 54000
$ws.Range("H18").Value2 = 61300
$ws.Range("I18").Value2 = 114500
$ws.Range("J18").Value2 = 129300
$ws.Range("K18").Value2 = 48700
$ws.Range("L18").Value2 = 58100
$ws.Range("M18").Value2 = 119300
$ws.Range("D20").Value2 = 0
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 0
$ws.Range("G20").Value2 = 0
$ws.Range("H20").Value2 = 0
$ws.Range("I20").Value2 = 0
$ws.Range("J20").Value2 = 0
$ws.Range("K20").Value2 = 0
$ws.Range("L20").Value2 = 0
$ws.Range("M20").Value2 = 0
$ws.Range("D21").Value2 = 64100
$ws.Range("E21").Value2 = 127800
$ws.Range("F21").Value2 = 142800
$ws.Range("G21").Value2 = 59500
$ws.Range("H21").Value2 = 66900
$ws.Range("I21").Value2 = 120000
$ws.Range("J21").Value2 = 134900
$ws.Range("K21").Value2 = 54100
$ws.Range("L21").Value2 = 63100
$ws.Range("M21").Value2 = 124300
$ws.Range("D22").Value2 = 400
$ws.Range("E22").Value2 = 1000
$ws.Range("F22").Value2 = 800
$ws.Range("G22").Value2 = 600
$ws.Range("H22").Value2 = 1300
$ws.Range("I22").Value2 = 2100
$ws.Range("J22").Value2 = 1600
$ws.Range("K22").Value2 = 1300
$ws.Range("L22").Value2 = 700
$ws.Range("M22").Value2 = 1000
$ws.Range("D23").Value2 = 58100
$ws.Range("E23").Value2 = 121300
$ws.Range("F23").Value2 = 136600
$ws.Range("G23").Value2 = 53400
$ws.Range("H23").Value2 = 60000
$ws.Range("I23").Value2 = 112300
$ws.Range("J23").Value2 = 127700
$ws.Range("K23").Value2 = 47400
$ws.Range("L23").Value2 = 57400
$ws.Range("M23").Value2 = 118300
$ws.Range("D24").Value2 = 11000
$ws.Range("E24").Value2 = 24400
$ws.Range("F24").Value2 = 28300
$ws.Range("G24").Value2 = 11000
$ws.Range("H24").Value2 = 17300
$ws.Range("I24").Value2 = 32300
$ws.Range("J24").Value2 = 36900
$ws.Range("K24").Value2 = 13700
$ws.Range("L24").Value2 = 17500
$ws.Range("M24").Value2 = 37800
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("F25").Value2 = 0
$ws.Range("G25").Value2 = 0
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 0
$ws.Range("L25").Value2 = 0
$ws.Range("M25").Value2 = 0
$ws.Range("D26").Value2 = 47200
$ws.Range("E26").Value2 = 96900
$ws.Range("F26").Value2 = 108300
$ws.Range("G26").Value2 = 42400
$ws.Range("H26").Value2 = 42700
$ws.Range("I26").Value2 = 80000
$ws.Range("J26").Value2 = 90800
$ws.Range("K26").Value2 = 33800
$ws.Range("L26").Value2 = 39900
$ws.Range("M26").Value2 = 80500
$ws.Range("D27").Value2 = 34600
$ws.Range("E27").Value2 = 72700
$ws.Range("F27").Value2 = 82600
$ws.Range("G27").Value2 = 31500
$ws.Range("H27").Value2 = 29700
$ws.Range("I27").Value2 = 59600
$ws.Range("J27").Value2 = 67600
$ws.Range("K27").Value2 = 24000
$ws.Range("L27").Value2 = 27100
$ws.Range("M27").Value2 = 58000
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 0
$ws.Range("G28").Value2 = 0
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 0
$ws.Range("J28").Value2 = 0
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = 0
$ws.Range("M28").Value2 = 0
$ws.Range("D29").Value2 = 1800
$ws.Range("E29").Value2 = "NA"
$ws.Range("F29").Value2 = "NA"
$ws.Range("G29").Value2 = "NA"
$ws.Range("H29").Value2 = 10000
$ws.Range("I29").Value2 = "NA"
$ws.Range("J29").Value2 = "NA"
$ws.Range("K29").Value2 = "NA"
$ws.Range("L29").Value2 = "NA"
$ws.Range("M29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("F30").Value2 = 0
$ws.Range("G30").Value2 = 0
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("M30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("F31").Value2 = 0
$ws.Range("G31").Value2 = 0
$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("E32").Value2 = 0
$ws.Range("F32").Value2 = 0
$ws.Range("G32").Value2 = 0
$ws.Range("H32").Value2 = 0
$ws.Range("I32").Value2 = 0
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 0
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = 0
$ws.Range("D33").Value2 = 36400
$ws.Range("E33").Value2 = 72700
$ws.Range("F33").Value2 = 82600
$ws.Range("G33").Value2 = 31500
$ws.Range("H33").Value2 = 39700
$ws.Range("I33").Value2 = 59600
$ws.Range("J33").Value2 = 67600
$ws.Range("K33").Value2 = 24000
$ws.Range("L33").Value2 = 27100
$ws.Range("M33").Value2 = 58000
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("F34").Value2 = 0
$ws.Range("G34").Value2 = 0
$ws.Range("H34").Value2 = 0
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 0
$ws.Range("M34").Value2 = 0
$ws.Range("D35").Value2 = 36400
$ws.Range("E35").Value2 = 72700
$ws.Range("F35").Value2 = 82600
$ws.Range("G35").Value2 = 31500
$ws.Range("H35").Value2 = 39700
$ws.Range("I35").Value2 = 59600
$ws.Range("J35").Value2 = 67600
$ws.Range("K35").Value2 = 24000
$ws.Range("L35").Value2 = 27100
$ws.Range("M35").Value2 = 58000
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("F38").Value2 = 43281
$ws.Range("G38").Value2 = 43190
$ws.Range("H38").Value2 = 43100
$ws.Range("I38").Value2 = 43008
$ws.Range("J38").Value2 = 42916
$ws.Range("K38").Value2 = 42825
$ws.Range("L38").Value2 = 42735
$ws.Range("M38").Value2 = 42643
$ws.Range("D41").Value2 = 82900
$ws.Range("E41").Value2 = 67600
$ws.Range("F41").Value2 = 72000
$ws.Range("G41").Value2 = 58100
$ws.Range("H41").Value2 = 80500
$ws.Range("I41").Value2 = 66700
$ws.Range("J41").Value2 = 51000
$ws.Range("K41").Value2 = 47400
$ws.Range("L41").Value2 = 56000
$ws.Range("M41").Value2 = 35900
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("F42").Value2 = 0
$ws.Range("G42").Value2 = 0
$ws.Range("H42").Value2 = 0
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = 0
$ws.Range("K42").Value2 = 0
$ws.Range("L42").Value2 = 0
$ws.Range("M42").Value2 = 0
$ws.Range("D43").Value2 = 501900
$ws.Range("E43").Value2 = 602800
$ws.Range("F43").Value2 = 659900
$ws.Range("G43").Value2 = 486700
$ws.Range("H43").Value2 = 478100
$ws.Range("I43").Value2 = 568500
$ws.Range("J43").Value2 = 613200
$ws.Range("K43").Value2 = 481600
$ws.Range("L43").Value2 = 476000
$ws.Range("M43").Value2 = 548000
$ws.Range("D44").Value2 = 837100
$ws.Range("E44").Value2 = 810900
$ws.Range("F44").Value2 = 872900
$ws.Range("G44").Value2 = 810000
$ws.Range("H44").Value2 = 761300
$ws.Range("I44").Value2 = 786100
$ws.Range("J44").Value2 = 770000
$ws.Range("K44").Value2 = 751500
$ws.Range("L44").Value2 = 685000
$ws.Range("M44").Value2 = 739700
$ws.Range("D45").Value2 = 19900
$ws.Range("E45").Value2 = 20600
$ws.Range("F45").Value2 = 16700
$ws.Range("G45").Value2 = 17000
$ws.Range("H45").Value2 = 17500
$ws.Range("I45").Value2 = 17800
$ws.Range("J45").Value2 = 17500
$ws.Range("K45").Value2 = 20000
$ws.Range("L45").Value2 = 23200
$ws.Range("M45").Value2 = 25400
$ws.Range("D46").Value2 = 1441800
$ws.Range("E46").Value2 = 1501800
$ws.Range("F46").Value2 = 1621500
$ws.Range("G46").Value2 = 1371700
$ws.Range("H46").Value2 = 1337400
$ws.Range("I46").Value2 = 1438900
$ws.Range("J46").Value2 = 1451700
$ws.Range("K46").Value2 = 1300500
$ws.Range("L46").Value2 = 1240200
$ws.Range("M46").Value2 = 1349000
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("F47").Value2 = 0
$ws.Range("G47").Value2 = 0
$ws.Range("H47").Value2 = 0
$ws.Range("I47").Value2 = 0
$ws.Range("J47").Value2 = 0
$ws.Range("K47").Value2 = 0
$ws.Range("L47").Value2 = 0
$ws.Range("M47").Value2 = 0
$ws.Range("D48").Value2 = 91000
$ws.Range("E48").Value2 = 91300
$ws.Range("F48").Value2 = 91300
$ws.Range("G48").Value2 = 90200
$ws.Range("H48").Value2 = 91200
$ws.Range("I48").Value2 = 91500
$ws.Range("J48").Value2 = 92300
$ws.Range("K48").Value2 = 90500
$ws.Range("L48").Value2 = 90500
$ws.Range("M48").Value2 = 59700
$ws.Range("D49").Value2 = 539800
$ws.Range("E49").Value2 = 550900
$ws.Range("F49").Value2 = 533500
$ws.Range("G49").Value2 = 538100
$ws.Range("H49").Value2 = 543800
$ws.Range("I49").Value2 = 546000
$ws.Range("J49").Value2 = 541000
$ws.Range("K49").Value2 = 538500
$ws.Range("L49").Value2 = 538300
$ws.Range("M49").Value2 = 543300
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("F50").Value2 = 0
$ws.Range("G50").Value2 = 0
$ws.Range("H50").Value2 = 0
$ws.Range("I50").Value2 = 0
$ws.Range("J50").Value2 = 0
$ws.Range("K50").Value2 = 0
$ws.Range("L50").Value2 = 0
$ws.Range("M50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("F51").Value2 = 0
$ws.Range("G51").Value2 = 0
$ws.Range("H51").Value2 = 0
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("K51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("M51").Value2 = 0
$ws.Range("D52").Value2 = 88300
$ws.Range("E52").Value2 = 86700
$ws.Range("F52").Value2 = 79300
$ws.Range("G52").Value2 = 76200
$ws.Range("H52").Value2 = 74500
$ws.Range("I52").Value2 = 71800
$ws.Range("J52").Value2 = 69400
$ws.Range("K52").Value2 = 5800
$ws.Range("L52").Value2 = 5700
$ws.Range("M52").Value2 = 5600
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("F53").Value2 = 0
$ws.Range("G53").Value2 = 0
$ws.Range("H53").Value2 = 0
$ws.Range("I53").Value2 = 0
$ws.Range("J53").Value2 = 0
$ws.Range("K53").Value2 = 0
$ws.Range("L53").Value2 = 0
$ws.Range("M53").Value2 = 0
$ws.Range("D54").Value2 = 2161000
$ws.Range("E54").Value2 = 2230700
$ws.Range("F54").Value2 = 2325600
$ws.Range("G54").Value2 = 2076300
$ws.Range("H54").Value2 = 2046900
$ws.Range("I54").Value2 = 2148200
$ws.Range("J54").Value2 = 2154400
$ws.Range("K54").Value2 = 1935300
$ws.Range("L54").Value2 = 1874600
$ws.Range("M54").Value2 = 1957600
$ws.Range("D57").Value2 = 200200
$ws.Range("E57").Value2 = 234500
$ws.Range("F57").Value2 = 358100
$ws.Range("G57").Value2 = 245400
$ws.Range("H57").Value2 = 230500
$ws.Range("I57").Value2 = 296300
$ws.Range("J57").Value2 = 286400
$ws.Range("K57").Value2 = 259700
$ws.Range("L57").Value2 = 185500
$ws.Range("M57").Value2 = 229700
$ws.Range("D58").Value2 = 200
$ws.Range("E58").Value2 = 116600
$ws.Range("F58").Value2 = 1700
$ws.Range("G58").Value2 = 200
$ws.Range("H58").Value2 = 200
$ws.Range("I58").Value2 = 200
$ws.Range("J58").Value2 = 1900
$ws.Range("K58").Value2 = 1800
$ws.Range("L58").Value2 = 200
$ws.Range("M58").Value2 = 200
$ws.Range("D59").Value2 = 157100
$ws.Range("E59").Value2 = 157800
$ws.Range("F59").Value2 = 163100
$ws.Range("G59").Value2 = 132100
$ws.Range("H59").Value2 = 185800
$ws.Range("I59").Value2 = 163800
$ws.Range("J59").Value2 = 136200
$ws.Range("K59").Value2 = 117400
$ws.Range("L59").Value2 = 129200
$ws.Range("M59").Value2 = 128200
$ws.Range("D60").Value2 = 357600
$ws.Range("E60").Value2 = 508800
$ws.Range("F60").Value2 = 522900
$ws.Range("G60").Value2 = 377700
$ws.Range("H60").Value2 = 416500
$ws.Range("I60").Value2 = 460400
$ws.Range("J60").Value2 = 424400
$ws.Range("K60").Value2 = 378900
$ws.Range("L60").Value2 = 314900
$ws.Range("M60").Value2 = 358100
$ws.Range("D61").Value2 = 135800
$ws.Range("E61").Value2 = 200
$ws.Range("F61").Value2 = 141800
$ws.Range("G61").Value2 = 91200
$ws.Range("H61").Value2 = 22100
$ws.Range("I61").Value2 = 285000
$ws.Range("J61").Value2 = 379500
$ws.Range("K61").Value2 = 280600
$ws.Range("L61").Value2 = 235600
$ws.Range("M61").Value2 = 219800
$ws.Range("D62").Value2 = 66000
$ws.Range("E62").Value2 = 61200
$ws.Range("F62").Value2 = 58900
$ws.Range("G62").Value2 = 57700
$ws.Range("H62").Value2 = 57300
$ws.Range("I62").Value2 = 70700
$ws.Range("J62").Value2 = 70100
$ws.Range("K62").Value2 = 67600
$ws.Range("L62").Value2 = 72400
$ws.Range("M62").Value2 = 73300
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("F63").Value2 = 0
$ws.Range("G63").Value2 = 0
$ws.Range("H63").Value2 = 0
$ws.Range("I63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("K63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("M63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("F64").Value2 = 0
$ws.Range("G64").Value2 = 0
$ws.Range("H64").Value2 = 0
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 0
$ws.Range("K64").Value2 = 0
$ws.Range("L64").Value2 = 0
$ws.Range("M64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("F65").Value2 = 0
$ws.Range("G65").Value2 = 0
$ws.Range("H65").Value2 = 0
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = 0
$ws.Range("D66").Value2 = 813200
$ws.Range("E66").Value2 = 863500
$ws.Range("F66").Value2 = 997200
$ws.Range("G66").Value2 = 783600
$ws.Range("H66").Value2 = 748900
$ws.Range("I66").Value2 = 1096100
$ws.Range("J66").Value2 = 1135800
$ws.Range("K66").Value2 = 956900
$ws.Range("L66").Value2 = 868800
$ws.Range("M66").Value2 = 918600
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("F68").Value2 = 0
$ws.Range("G68").Value2 = 0
$ws.Range("H68").Value2 = 0
$ws.Range("I68").Value2 = 0
$ws.Range("J68").Value2 = 0
$ws.Range("K68").Value2 = 0
$ws.Range("L68").Value2 = 0
$ws.Range("M68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("F69").Value2 = 0
$ws.Range("G69").Value2 = 0
$ws.Range("H69").Value2 = 0
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 0
$ws.Range("M69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("F70").Value2 = 0
$ws.Range("G70").Value2 = 0
$ws.Range("H70").Value2 = 0
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("K70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("M70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("F71").Value2 = 0
$ws.Range("G71").Value2 = 0
$ws.Range("H71").Value2 = 0
$ws.Range("I71").Value2 = 0
$ws.Range("J71").Value2 = 0
$ws.Range("K71").Value2 = 0
$ws.Range("L71").Value2 = 0
$ws.Range("M71").Value2 = 0
$ws.Range("D72").Value2 = 628000
$ws.Range("E72").Value2 = 642600
$ws.Range("F72").Value2 = 617700
$ws.Range("G72").Value2 = 581900
$ws.Range("H72").Value2 = 594600
$ws.Range("I72").Value2 = 596000
$ws.Range("J72").Value2 = 575600
$ws.Range("K72").Value2 = 539300
$ws.Range("L72").Value2 = 550500
$ws.Range("M72").Value2 = 558200
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("F73").Value2 = 0
$ws.Range("G73").Value2 = 0
$ws.Range("H73").Value2 = 0
$ws.Range("I73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("K73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("M73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("F74").Value2 = 0
$ws.Range("G74").Value2 = 0
$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 0
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("L74").Value2 = 0
$ws.Range("M74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("F75").Value2 = 0
$ws.Range("G75").Value2 = 0
$ws.Range("H75").Value2 = 0
$ws.Range("I75").Value2 = 0
$ws.Range("J75").Value2 = 0
$ws.Range("K75").Value2 = 0
$ws.Range("L75").Value2 = 0
$ws.Range("M75").Value2 = 0
$ws.Range("D76").Value2 = 1347800
$ws.Range("E76").Value2 = 1367200
$ws.Range("F76").Value2 = 1328400
$ws.Range("G76").Value2 = 1292700
$ws.Range("H76").Value2 = 1298000
$ws.Range("I76").Value2 = 1052100
$ws.Range("J76").Value2 = 1018700
$ws.Range("K76").Value2 = 978300
$ws.Range("L76").Value2 = 1005800
$ws.Range("M76").Value2 = 1039000
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("F77").Value2 = 0
$ws.Range("G77").Value2 = 0
$ws.Range("H77").Value2 = 0
$ws.Range("I77").Value2 = 0
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("L77").Value2 = 0
$ws.Range("M77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("F80").Value2 = 43281
$ws.Range("G80").Value2 = 43190
$ws.Range("H80").Value2 = 43100
$ws.Range("I80").Value2 = 43008
$ws.Range("J80").Value2 = 42916
$ws.Range("K80").Value2 = 42825
$ws.Range("L80").Value2 = 42735
$ws.Range("M80").Value2 = 42643
$ws.Range("D81").Value2 = 36400
$ws.Range("E81").Value2 = 72700
$ws.Range("F81").Value2 = 82600
$ws.Range("G81").Value2 = 31500
$ws.Range("H81").Value2 = 39700
$ws.Range("I81").Value2 = 59600
$ws.Range("J81").Value2 = 67600
$ws.Range("K81").Value2 = 24000
$ws.Range("L81").Value2 = 27100
$ws.Range("M81").Value2 = 58000
$ws.Range("D83").Value2 = 5600
$ws.Range("E83").Value2 = 5500
$ws.Range("F83").Value2 = 5500
$ws.Range("G83").Value2 = 5500
$ws.Range("H83").Value2 = 5500
$ws.Range("I83").Value2 = 5600
$ws.Range("J83").Value2 = 5600
$ws.Range("K83").Value2 = 5400
$ws.Range("L83").Value2 = 5000
$ws.Range("M83").Value2 = 5000
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("F84").Value2 = 0
$ws.Range("G84").Value2 = 0
$ws.Range("H84").Value2 = 0
$ws.Range("I84").Value2 = 0
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = 0
$ws.Range("L84").Value2 = 0
$ws.Range("M84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("F85").Value2 = 0
$ws.Range("G85").Value2 = 0
$ws.Range("H85").Value2 = 0
$ws.Range("I85").Value2 = 0
$ws.Range("J85").Value2 = 0
$ws.Range("K85").Value2 = 0
$ws.Range("L85").Value2 = 0
$ws.Range("M85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("F86").Value2 = 0
$ws.Range("G86").Value2 = 0
$ws.Range("H86").Value2 = 0
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 0
$ws.Range("L86").Value2 = 0
$ws.Range("M86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("F87").Value2 = 0
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = 0
$ws.Range("I87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("K87").Value2 = 0
$ws.Range("L87").Value2 = 0
$ws.Range("M87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("F88").Value2 = 0
$ws.Range("G88").Value2 = 0
$ws.Range("H88").Value2 = 0
$ws.Range("I88").Value2 = 0
$ws.Range("J88").Value2 = 0
$ws.Range("K88").Value2 = 0
$ws.Range("L88").Value2 = 0
$ws.Range("M88").Value2 = 0
$ws.Range("D89").Value2 = 100500
$ws.Range("E89").Value2 = 89300
$ws.Range("F89").Value2 = 22400
$ws.Range("G89").Value2 = -41600
$ws.Range("H89").Value2 = 117100
$ws.Range("I89").Value2 = 154400
$ws.Range("J89").Value2 = 700
$ws.Range("K89").Value2 = 34300
$ws.Range("L89").Value2 = 135400
$ws.Range("M89").Value2 = 104400
$ws.Range("D91").Value2 = -4300
$ws.Range("E91").Value2 = -4100
$ws.Range("F91").Value2 = -5300
$ws.Range("G91").Value2 = -3500
$ws.Range("H91").Value2 = -4000
$ws.Range("I91").Value2 = -3500
$ws.Range("J91").Value2 = -6200
$ws.Range("K91").Value2 = -4100
$ws.Range("L91").Value2 = -78200
$ws.Range("M91").Value2 = -3400
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("F92").Value2 = 0
$ws.Range("G92").Value2 = 0
$ws.Range("H92").Value2 = 0
$ws.Range("I92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("K92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("M92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("F93").Value2 = 0
$ws.Range("G93").Value2 = 0
$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value2 = 0
$ws.Range("J93").Value2 = 0
$ws.Range("K93").Value2 = 0
$ws.Range("L93").Value2 = 0
$ws.Range("M93").Value2 = 0
$ws.Range("D94").Value2 = -4000
$ws.Range("E94").Value2 = -13600
$ws.Range("F94").Value2 = -5300
$ws.Range("G94").Value2 = -3400
$ws.Range("H94").Value2 = -4000
$ws.Range("I94").Value2 = -3500
$ws.Range("J94").Value2 = -69700
$ws.Range("K94").Value2 = -4100
$ws.Range("L94").Value2 = -34500
$ws.Range("M94").Value2 = -3300
$ws.Range("D96").Value2 = -54300
$ws.Range("E96").Value2 = -54200
$ws.Range("F96").Value2 = -54200
$ws.Range("G96").Value2 = -46600
$ws.Range("H96").Value2 = -44700
$ws.Range("I96").Value2 = -44600
$ws.Range("J96").Value2 = -37500
$ws.Range("K96").Value2 = -37400
$ws.Range("L96").Value2 = -37300
$ws.Range("M96").Value2 = -30100
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("F97").Value2 = 0
$ws.Range("G97").Value2 = 0
$ws.Range("H97").Value2 = 0
$ws.Range("I97").Value2 = 0
$ws.Range("J97").Value2 = 0
$ws.Range("K97").Value2 = 0
$ws.Range("L97").Value2 = 0
$ws.Range("M97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("F98").Value2 = 0
$ws.Range("G98").Value2 = 0
$ws.Range("H98").Value2 = 0
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 0
$ws.Range("K98").Value2 = 0
$ws.Range("L98").Value2 = 0
$ws.Range("M98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("F99").Value2 = 0
$ws.Range("G99").Value2 = 0
$ws.Range("H99").Value2 = 0
$ws.Range("I99").Value2 = 0
$ws.Range("J99").Value2 = 0
$ws.Range("K99").Value2 = 0
$ws.Range("L99").Value2 = 0
$ws.Range("M99").Value2 = 0
$ws.Range("D100").Value2 = -79800
$ws.Range("E100").Value2 = -80500
$ws.Range("F100").Value2 = -2600
$ws.Range("G100").Value2 = 23300
$ws.Range("H100").Value2 = -103900
$ws.Range("I100").Value2 = -134200
$ws.Range("J100").Value2 = 74800
$ws.Range("K100").Value2 = -38900
$ws.Range("L100").Value2 = -76500
$ws.Range("M100").Value2 = -95700
$ws.Range("D101").Value2 = -1400
$ws.Range("E101").Value2 = 400
$ws.Range("F101").Value2 = -600
$ws.Range("G101").Value2 = -700
$ws.Range("H101").Value2 = -100
$ws.Range("I101").Value2 = 800
$ws.Range("J101").Value2 = 600
$ws.Range("K101").Value2 = 200
$ws.Range("L101").Value2 = -300
$ws.Range("M101").Value2 = 0
$ws.Range("D102").Value2 = 15300
$ws.Range("E102").Value2 = -4400
$ws.Range("F102").Value2 = 13900
$ws.Range("G102").Value2 = -22400
$ws.Range("H102").Value2 = 13800
$ws.Range("I102").Value2 = 15600
$ws.Range("J102").Value2 = 3600
$ws.Range("K102").Value2 = -8600
$ws.Range("L102").Value2 = 20100
$ws.Range("M102").Value2 = 5400

# Column widths shift slightly (minor re-fit of "best fit" widths); set explicit
# widths close to the saved workbook's values.
$ws.Columns("A").ColumnWidth = 6.166666666666667
$ws.Columns("B").ColumnWidth = 26.053385416666668
$ws.Columns("C").ColumnWidth = 68.27604166666667
$ws.Range("D1:E1").ColumnWidth = 13.830729166666666
$ws.Columns("F").ColumnWidth = 13.608072916666666
$ws.Columns("G").ColumnWidth = 14.053385416666666
$ws.Range("H1:I1").ColumnWidth = 13.830729166666666
$ws.Columns("J").ColumnWidth = 13.608072916666666
$ws.Columns("K").ColumnWidth = 14.053385416666666
$ws.Range("L1:M1").ColumnWidth = 13.830729166666666

